# Update column G ("K") with regenerated values (s_vals), replacing the
# previous Strike# based figures, per: "regen save_data to use K instead of
# Strike#, regen std/mean, calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(0,0,0,0,0,2,0,0,2,1,0,1,1,1,1,1,0,1,2,3,6,5,4,4,2,3,3,0,4,4,7,7,5,4,6,6,3,4,2,6,5,4,3,5,0)

for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
